$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 831.5
$ws.Range("I19").Value = 516.6667
$ws.Range("J19").Value = 1020.4
$ws.Range("K19").Value = 516.6667
$ws.Range("L19").Value = 1020.4
$ws.Range("M19").Value = -341.6667
$ws.Range("N19").Value = -1370.4
$ws.Range("H58").Value = 1125.6666
$ws.Range("I58").Value = 532.6923
$ws.Range("J58").Value = 4980
$ws.Range("K58").Value = 1598.0769
$ws.Range("L58").Value = 14940
$ws.Range("M58").Value = -1448.0769
$ws.Range("N58").Value = -15240
$ws.Range("H61").Value = 113.666664
$ws.Range("I61").Value = 113.666664
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 340.999992
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -168.999992
$ws.Range("H64").Value = 4785.7144
$ws.Range("I64").Value = 3850
$ws.Range("J64").Value = 5160
$ws.Range("K64").Value = 3850
$ws.Range("L64").Value = 5160
$ws.Range("M64").Value = -3602
$ws.Range("N64").Value = -5656
$ws.Range("H67").Value = 4785.7144
$ws.Range("I67").Value = 3850
$ws.Range("J67").Value = 5160
$ws.Range("K67").Value = 3850
$ws.Range("L67").Value = 5160
$ws.Range("M67").Value = -2992
$ws.Range("N67").Value = -6876
$ws.Range("H70").Value = 4998.6
$ws.Range("I70").Value = 1495
$ws.Range("J70").Value = 7334.3335
$ws.Range("K70").Value = 4485
$ws.Range("L70").Value = 22003.0005
$ws.Range("M70").Value = -4215
$ws.Range("N70").Value = -22543.0005
$ws.Range("H73").Value = 4998.6
$ws.Range("I73").Value = 1495
$ws.Range("J73").Value = 7334.3335
$ws.Range("K73").Value = 4485
$ws.Range("L73").Value = 22003.0005
$ws.Range("M73").Value = -3549
$ws.Range("N73").Value = -23875.0005
$ws.Range("H86").Value = 1340
$ws.Range("I86").Value = 1333.5
$ws.Range("J86").Value = 1349.75
$ws.Range("K86").Value = 1333.5
$ws.Range("L86").Value = 1349.75
$ws.Range("M86").Value = -210.5
$ws.Range("N86").Value = -3595.75
$ws.Range("H89").Value = 1340
$ws.Range("I89").Value = 1333.5
$ws.Range("J89").Value = 1349.75
$ws.Range("K89").Value = 6667.5
$ws.Range("L89").Value = 6748.75
$ws.Range("M89").Value = -1051.5
$ws.Range("N89").Value = -17980.75
$ws.Range("H94").Value = 2415.9
$ws.Range("I94").Value = 2464.3333
$ws.Range("J94").Value = 1980
$ws.Range("K94").Value = 2464.3333
$ws.Range("L94").Value = 1980
$ws.Range("M94").Value = -2013.3333
$ws.Range("N94").Value = -2882
$ws.Range("H131").Value = 1327
$ws.Range("I131").Value = 696.6667
$ws.Range("J131").Value = 7000
$ws.Range("K131").Value = 2090.0001
$ws.Range("L131").Value = 21000
$ws.Range("M131").Value = 2949.9999
$ws.Range("N131").Value = -31080
$ws.Range("H135").Value = 18481.068
$ws.Range("I135").Value = 22505.436
$ws.Range("J135").Value = 3054.3333
$ws.Range("K135").Value = 202548.924
$ws.Range("L135").Value = 27488.9997
$ws.Range("M135").Value = -200013.924
$ws.Range("N135").Value = -32558.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 9300
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 9300
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 9300
$ws.Range("N43").Value = -9926
$ws.Range("H97").Value = 3125866
$ws.Range("I97").Value = 5209151.5
$ws.Range("J97").Value = 937.5
$ws.Range("K97").Value = 5209151.5
$ws.Range("L97").Value = 937.5
$ws.Range("M97").Value = -5208655.5
$ws.Range("N97").Value = -1929.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 520.6923
$ws.Range("I94").Value = 323.625
$ws.Range("J94").Value = 836
$ws.Range("K94").Value = 323.625
$ws.Range("L94").Value = 836
$ws.Range("M94").Value = 127.375
$ws.Range("N94").Value = -1738
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H134").Value = 1926.6888
$ws.Range("I134").Value = 1018.86664
$ws.Range("J134").Value = 3742.3333
$ws.Range("K134").Value = 3056.59992
$ws.Range("L134").Value = 11226.9999
$ws.Range("M134").Value = -521.5999199999997
$ws.Range("N134").Value = -16296.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3704.4517
$ws.Range("I31").Value = 1376.5834
$ws.Range("J31").Value = 11685.714
$ws.Range("K31").Value = 1376.5834
$ws.Range("L31").Value = 11685.714
$ws.Range("M31").Value = -1081.5834
$ws.Range("N31").Value = -12275.714
$ws.Range("H34").Value = 3704.4517
$ws.Range("I34").Value = 1376.5834
$ws.Range("J34").Value = 11685.714
$ws.Range("K34").Value = 1376.5834
$ws.Range("L34").Value = 11685.714
$ws.Range("M34").Value = -1174.5834
$ws.Range("N34").Value = -12089.714
$ws.Range("H98").Value = 54935
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 54935
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 54935
$ws.Range("N98").Value = -59427
$ws.Range("H99").Value = 8510.277
$ws.Range("I99").Value = 8629.615
$ws.Range("J99").Value = 8200
$ws.Range("K99").Value = 8629.615
$ws.Range("L99").Value = 8200
$ws.Range("M99").Value = -7131.615
$ws.Range("N99").Value = -11196
$ws.Range("H126").Value = 8510.277
$ws.Range("I126").Value = 8629.615
$ws.Range("J126").Value = 8200
$ws.Range("K126").Value = 25888.845
$ws.Range("L126").Value = 24600
$ws.Range("M126").Value = -23418.845
$ws.Range("N126").Value = -29540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 211.41667
$ws.Range("I2").Value = 248.55556
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1491.33336
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -1378.33336
$ws.Range("N2").Value = -826
$ws.Range("H37").Value = 94090.91
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 94090.91
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 282272.73
$ws.Range("N37").Value = -282496.73
$ws.Range("H101").Value = 5637.931
$ws.Range("I101").Value = 5000
$ws.Range("J101").Value = 5660.7144
$ws.Range("K101").Value = 15000
$ws.Range("L101").Value = 16982.1432
$ws.Range("M101").Value = -12566
$ws.Range("N101").Value = -21850.1432
$ws.Range("H131").Value = 872.3134
$ws.Range("I131").Value = 423.2857
$ws.Range("J131").Value = 990.9245
$ws.Range("K131").Value = 1269.8571
$ws.Range("L131").Value = 2972.7735
$ws.Range("M131").Value = 3770.1429
$ws.Range("N131").Value = -13052.7735
$ws.Range("H132").Value = 1331.2727
$ws.Range("I132").Value = 789
$ws.Range("J132").Value = 1982
$ws.Range("K132").Value = 7101
$ws.Range("L132").Value = 17838
$ws.Range("M132").Value = -4571
$ws.Range("N132").Value = -22898
$ws.Range("H134").Value = 1981
$ws.Range("I134").Value = 1452.909
$ws.Range("J134").Value = 7790
$ws.Range("K134").Value = 4358.727000000001
$ws.Range("L134").Value = 23370
$ws.Range("M134").Value = 711.2729999999992
$ws.Range("N134").Value = -33510

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50542.637
$ws.Range("I70").Value = 76195.57000000001
$ws.Range("J70").Value = 5650
$ws.Range("K70").Value = 76195.57000000001
$ws.Range("L70").Value = 5650
$ws.Range("M70").Value = -75925.57000000001
$ws.Range("N70").Value = -6190
$ws.Range("H73").Value = 50542.637
$ws.Range("I73").Value = 76195.57000000001
$ws.Range("J73").Value = 5650
$ws.Range("K73").Value = 76195.57000000001
$ws.Range("L73").Value = 5650
$ws.Range("M73").Value = -75259.57000000001
$ws.Range("N73").Value = -7522
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H109").Value = 10285
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 10285
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365
$ws.Range("H141").Value = 64500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 64500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 64500
$ws.Range("N141").Value = -74860

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1924.2693
$ws.Range("I7").Value = 1884.1364
$ws.Range("J7").Value = 2145
$ws.Range("K7").Value = 1884.1364
$ws.Range("L7").Value = 2145
$ws.Range("M7").Value = -1772.1364
$ws.Range("N7").Value = -2369
$ws.Range("H40").Value = 3067.1304
$ws.Range("I40").Value = 2857.7222
$ws.Range("J40").Value = 3821
$ws.Range("K40").Value = 2857.7222
$ws.Range("L40").Value = 3821
$ws.Range("M40").Value = -2721.7222
$ws.Range("N40").Value = -4093
$ws.Range("H126").Value = 1924.2693
$ws.Range("I126").Value = 1884.1364
$ws.Range("J126").Value = 2145
$ws.Range("K126").Value = 5652.4092
$ws.Range("L126").Value = 6435
$ws.Range("M126").Value = -3182.4092
$ws.Range("N126").Value = -11375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = -1280
$ws.Range("H81").Value = 69071.39999999999
$ws.Range("I81").Value = 1133
$ws.Range("J81").Value = 79523.46000000001
$ws.Range("K81").Value = 2266
$ws.Range("L81").Value = 159046.92
$ws.Range("M81").Value = -1205
$ws.Range("N81").Value = -161168.92
$ws.Range("H84").Value = 69071.39999999999
$ws.Range("I84").Value = 1133
$ws.Range("J84").Value = 79523.46000000001
$ws.Range("K84").Value = 11330
$ws.Range("L84").Value = 795234.6000000001
$ws.Range("M84").Value = -6026
$ws.Range("N84").Value = -805842.6000000001
$ws.Range("H96").Value = 1900.8
$ws.Range("I96").Value = 1133.3334
$ws.Range("J96").Value = 3052
$ws.Range("K96").Value = 1133.3334
$ws.Range("L96").Value = 3052
$ws.Range("M96").Value = 239.6666
$ws.Range("N96").Value = -5798
$ws.Range("H126").Value = 937.7273
$ws.Range("I126").Value = 771.5789
$ws.Range("J126").Value = 1990
$ws.Range("K126").Value = 2314.7367
$ws.Range("L126").Value = 5970
$ws.Range("M126").Value = 155.2633000000001
$ws.Range("N126").Value = -10910

